$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.062193222105793
$ws.Cells.Item(2, 4).Value = 1.071922342791158
$ws.Cells.Item(2, 5).Value = 1.068113397621157
$ws.Cells.Item(2, 6).Value = 1.079679277943972
$ws.Cells.Item(2, 9).Value = 1.043509016713692
$ws.Cells.Item(2, 10).Value = 1.067165007485296
$ws.Cells.Item(2, 11).Value = 1.074618227955066
$ws.Cells.Item(2, 12).Value = 1.070819449833063
$ws.Cells.Item(2, 13).Value = 1.082354697468612
$ws.Cells.Item(2, 14).Value = 1.026012963628157

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.063436798694366
$ws.Cells.Item(3, 4).Value = 1.073092899128294
$ws.Cells.Item(3, 5).Value = 1.06921986804821
$ws.Cells.Item(3, 6).Value = 1.080868974498445
$ws.Cells.Item(3, 9).Value = 1.043737231411955
$ws.Cells.Item(3, 10).Value = 1.068061521644544
$ws.Cells.Item(3, 11).Value = 1.07560455704297
$ws.Cells.Item(3, 12).Value = 1.071741113415765
$ws.Cells.Item(3, 13).Value = 1.083361606290908
$ws.Cells.Item(3, 14).Value = 1.026319030555606

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.064241303195013
$ws.Cells.Item(4, 4).Value = 1.073850434440089
$ws.Cells.Item(4, 5).Value = 1.069935963853285
$ws.Cells.Item(4, 6).Value = 1.081638961030051
$ws.Cells.Item(4, 9).Value = 1.0438833497747
$ws.Cells.Item(4, 10).Value = 1.068640915186397
$ws.Cells.Item(4, 11).Value = 1.076242302801203
$ws.Cells.Item(4, 12).Value = 1.072337032408312
$ws.Cells.Item(4, 13).Value = 1.084012736314587
$ws.Cells.Item(4, 14).Value = 1.026516655769923

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.064579477498714
$ws.Cells.Item(5, 4).Value = 1.074168929181758
$ws.Cells.Item(5, 5).Value = 1.070237044031527
$ws.Cells.Item(5, 6).Value = 1.081962706088641
$ws.Cells.Item(5, 9).Value = 1.043944406660332
$ws.Cells.Item(5, 10).Value = 1.068884323089063
$ws.Cells.Item(5, 11).Value = 1.076510298514093
$ws.Cells.Item(5, 12).Value = 1.072587447895828
$ws.Cells.Item(5, 13).Value = 1.084286374974266
$ws.Cells.Item(5, 14).Value = 1.026599636989721

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.064636256199522
$ws.Cells.Item(6, 4).Value = 1.074222407464143
$ws.Cells.Item(6, 5).Value = 1.070287598714039
$ws.Cells.Item(6, 6).Value = 1.082017066849433
$ws.Cells.Item(6, 9).Value = 1.043954636615126
$ws.Cells.Item(6, 10).Value = 1.068925182459999
$ws.Cells.Item(6, 11).Value = 1.076555289603486
$ws.Cells.Item(6, 12).Value = 1.07262948740434
$ws.Cells.Item(6, 13).Value = 1.084332314488448
$ws.Cells.Item(6, 14).Value = 1.026613564020533

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.064245822048929
$ws.Cells.Item(7, 4).Value = 1.073854690075158
$ws.Cells.Item(7, 5).Value = 1.069939986768383
$ws.Cells.Item(7, 6).Value = 1.08164328675797
$ws.Cells.Item(7, 9).Value = 1.043884167078359
$ws.Cells.Item(7, 10).Value = 1.068644168277964
$ws.Cells.Item(7, 11).Value = 1.076245884214073
$ws.Cells.Item(7, 12).Value = 1.072340378898773
$ws.Cells.Item(7, 13).Value = 1.084016393065164
$ws.Cells.Item(7, 14).Value = 1.026517764963301

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.062613531575504
$ws.Cells.Item(8, 4).Value = 1.072317916371206
$ws.Cells.Item(8, 5).Value = 1.068487307084083
$ws.Cells.Item(8, 6).Value = 1.080081306321536
$ws.Cells.Item(8, 9).Value = 1.043586464391374
$ws.Cells.Item(8, 10).Value = 1.067468136530141
$ws.Cells.Item(8, 11).Value = 1.074951661550978
$ws.Cells.Item(8, 12).Value = 1.071131026006331
$ws.Cells.Item(8, 13).Value = 1.082695072025326
$ws.Cells.Item(8, 14).Value = 1.0261164876155

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.05973581654328
$ws.Cells.Item(9, 4).Value = 1.069610688081723
$ws.Cells.Item(9, 5).Value = 1.065928488672691
$ws.Cells.Item(9, 6).Value = 1.077330168990712
$ws.Cells.Item(9, 9).Value = 1.043049978589529
$ws.Cells.Item(9, 10).Value = 1.06539032325106
$ws.Cells.Item(9, 11).Value = 1.072667385576384
$ws.Cells.Item(9, 12).Value = 1.068996430273102
$ws.Cells.Item(9, 13).Value = 1.080363565855598
$ws.Cells.Item(9, 14).Value = 1.025406152793433

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.057816266130723
$ws.Cells.Item(10, 4).Value = 1.06780630021494
$ws.Cells.Item(10, 5).Value = 1.064223203877678
$ws.Cells.Item(10, 6).Value = 1.075496852006688
$ws.Cells.Item(10, 9).Value = 1.042684307309748
$ws.Cells.Item(10, 10).Value = 1.064001348740356
$ws.Cells.Item(10, 11).Value = 1.071141979333286
$ws.Cells.Item(10, 12).Value = 1.067570904936525
$ws.Cells.Item(10, 13).Value = 1.078807024632128
$ws.Cells.Item(10, 14).Value = 1.024930404931282

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.056984796325936
$ws.Cells.Item(11, 4).Value = 1.067025063419381
$ws.Cells.Item(11, 5).Value = 1.06348492108546
$ws.Cells.Item(11, 6).Value = 1.074703170441971
$ws.Cells.Item(11, 9).Value = 1.042524062778671
$ws.Cells.Item(11, 10).Value = 1.063398997721227
$ws.Cells.Item(11, 11).Value = 1.070480838082569
$ws.Cells.Item(11, 12).Value = 1.066953037223678
$ws.Cells.Item(11, 13).Value = 1.078132485413961
$ws.Cells.Item(11, 14).Value = 1.024723876831946

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.05667590501555
$ws.Cells.Item(12, 4).Value = 1.066734886812602
$ws.Cells.Item(12, 5).Value = 1.063210706001156
$ws.Cells.Item(12, 6).Value = 1.074408383603045
$ws.Cells.Item(12, 9).Value = 1.04246425399379
$ws.Cells.Item(12, 10).Value = 1.063175118607264
$ws.Cells.Item(12, 11).Value = 1.070235164496992
$ws.Cells.Item(12, 12).Value = 1.066723441179799
$ws.Cells.Item(12, 13).Value = 1.077881847830332
$ws.Cells.Item(12, 14).Value = 1.024647083644058

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.05674216539126
$ws.Cells.Item(13, 4).Value = 1.066797130277033
$ws.Cells.Item(13, 5).Value = 1.063269525367389
$ws.Cells.Item(13, 6).Value = 1.074471615420417
$ws.Cells.Item(13, 9).Value = 1.042477096167456
$ws.Cells.Item(13, 10).Value = 1.06322314775744
$ws.Cells.Item(13, 11).Value = 1.070287866697953
$ws.Cells.Item(13, 12).Value = 1.066772694518077
$ws.Cells.Item(13, 13).Value = 1.07793561425489
$ws.Cells.Item(13, 14).Value = 1.024663559645713

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.056959264235323
$ws.Cells.Item(14, 4).Value = 1.067001077146424
$ws.Cells.Item(14, 5).Value = 1.06346225406317
$ws.Cells.Item(14, 6).Value = 1.074678802834508
$ws.Cells.Item(14, 9).Value = 1.04251912481949
$ws.Cells.Item(14, 10).Value = 1.063380494650783
$ws.Cells.Item(14, 11).Value = 1.070460532600834
$ws.Cells.Item(14, 12).Value = 1.066934060637698
$ws.Cells.Item(14, 13).Value = 1.07811176934931
$ws.Cells.Item(14, 14).Value = 1.024717530706342

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.057093019772632
$ws.Cells.Item(15, 4).Value = 1.067126736753197
$ws.Cells.Item(15, 5).Value = 1.063581002650074
$ws.Cells.Item(15, 6).Value = 1.074806460679819
$ws.Cells.Item(15, 9).Value = 1.042544982038988
$ws.Cells.Item(15, 10).Value = 1.0634774227837
$ws.Cells.Item(15, 11).Value = 1.070566904955432
$ws.Cells.Item(15, 12).Value = 1.067033471333968
$ws.Cells.Item(15, 13).Value = 1.078220293180997
$ws.Cells.Item(15, 14).Value = 1.024750773517046

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.057871441676586
$ws.Cells.Item(16, 4).Value = 1.067858149760014
$ws.Cells.Item(16, 5).Value = 1.06427220360522
$ws.Cells.Item(16, 6).Value = 1.075549529194638
$ws.Cells.Item(16, 9).Value = 1.042694902009941
$ws.Cells.Item(16, 10).Value = 1.064041305367007
$ws.Cells.Item(16, 11).Value = 1.071185843685345
$ws.Cells.Item(16, 12).Value = 1.067611897894264
$ws.Cells.Item(16, 13).Value = 1.078851779928941
$ws.Cells.Item(16, 14).Value = 1.024944100401021

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.058359645237542
$ws.Cells.Item(17, 4).Value = 1.068316964674331
$ws.Cells.Item(17, 5).Value = 1.064705805890391
$ws.Cells.Item(17, 6).Value = 1.076015677203876
$ws.Cells.Item(17, 9).Value = 1.042788432147282
$ws.Cells.Item(17, 10).Value = 1.064394767729502
$ws.Cells.Item(17, 11).Value = 1.07157391772925
$ws.Cells.Item(17, 12).Value = 1.06797456637545
$ws.Cells.Item(17, 13).Value = 1.079247747439629
$ws.Cells.Item(17, 14).Value = 1.025065228097258

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.058644378134124
$ws.Cells.Item(18, 4).Value = 1.068584591013292
$ws.Cells.Item(18, 5).Value = 1.064958730139384
$ws.Cells.Item(18, 6).Value = 1.076287588696732
$ws.Cells.Item(18, 9).Value = 1.042842802776956
$ws.Cells.Item(18, 10).Value = 1.064600848079536
$ws.Cells.Item(18, 11).Value = 1.071800214049502
$ws.Cells.Item(18, 12).Value = 1.068186046409287
$ws.Cells.Item(18, 13).Value = 1.079478655957791
$ws.Cells.Item(18, 14).Value = 1.025135829096958

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.058741460024637
$ws.Cells.Item(19, 4).Value = 1.068675845979952
$ws.Cells.Item(19, 5).Value = 1.06504497274668
$ws.Cells.Item(19, 6).Value = 1.076380306126417
$ws.Cells.Item(19, 9).Value = 1.042861310573224
$ws.Cells.Item(19, 10).Value = 1.064671101227479
$ws.Cells.Item(19, 11).Value = 1.071877364997988
$ws.Cells.Item(19, 12).Value = 1.06825814578082
$ws.Cells.Item(19, 13).Value = 1.079557380914826
$ws.Cells.Item(19, 14).Value = 1.025159893619832

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.058307268529705
$ws.Cells.Item(20, 4).Value = 1.068267737428477
$ws.Cells.Item(20, 5).Value = 1.064659283304555
$ws.Cells.Item(20, 6).Value = 1.075965662370679
$ws.Cells.Item(20, 9).Value = 1.04277841627541
$ws.Cells.Item(20, 10).Value = 1.064356853704828
$ws.Cells.Item(20, 11).Value = 1.071532287348427
$ws.Cells.Item(20, 12).Value = 1.067935661527002
$ws.Cells.Item(20, 13).Value = 1.079205269329354
$ws.Cells.Item(20, 14).Value = 1.025052237492957

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.056895335338369
$ws.Cells.Item(21, 4).Value = 1.066941019627838
$ws.Cells.Item(21, 5).Value = 1.063405499836056
$ws.Cells.Item(21, 6).Value = 1.074617790716239
$ws.Cells.Item(21, 9).Value = 1.042506756350231
$ws.Cells.Item(21, 10).Value = 1.063334163765275
$ws.Cells.Item(21, 11).Value = 1.07040968946687
$ws.Cells.Item(21, 12).Value = 1.066886544901524
$ws.Cells.Item(21, 13).Value = 1.078059898387234
$ws.Cells.Item(21, 14).Value = 1.024701639770864

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.056007327151375
$ws.Cells.Item(22, 4).Value = 1.066106913718068
$ws.Cells.Item(22, 5).Value = 1.062617288314339
$ws.Cells.Item(22, 6).Value = 1.073770455771328
$ws.Cells.Item(22, 9).Value = 1.042334293253392
$ws.Cells.Item(22, 10).Value = 1.062690352313032
$ws.Cells.Item(22, 11).Value = 1.06970330973818
$ws.Cells.Item(22, 12).Value = 1.066226387948093
$ws.Cells.Item(22, 13).Value = 1.077339273267337
$ws.Cells.Item(22, 14).Value = 1.024480745154267

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.05647810351958
$ws.Cells.Item(23, 4).Value = 1.0665490844895
$ws.Cells.Item(23, 5).Value = 1.063035125955556
$ws.Cells.Item(23, 6).Value = 1.074219632742476
$ws.Cells.Item(23, 9).Value = 1.04242587666311
$ws.Cells.Item(23, 10).Value = 1.063031725826193
$ws.Cells.Item(23, 11).Value = 1.070077828412759
$ws.Cells.Item(23, 12).Value = 1.066576400901065
$ws.Cells.Item(23, 13).Value = 1.077721336769537
$ws.Cells.Item(23, 14).Value = 1.024597889302158

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.058330935395921
$ws.Cells.Item(24, 4).Value = 1.068289981079259
$ws.Cells.Item(24, 5).Value = 1.064680304823515
$ws.Cells.Item(24, 6).Value = 1.075988261873928
$ws.Cells.Item(24, 9).Value = 1.042782942585098
$ws.Cells.Item(24, 10).Value = 1.064373985693376
$ws.Cells.Item(24, 11).Value = 1.071551098513199
$ws.Cells.Item(24, 12).Value = 1.067953241133757
$ws.Cells.Item(24, 13).Value = 1.07922446352164
$ws.Cells.Item(24, 14).Value = 1.025058107544672

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.060479954913092
$ws.Cells.Item(25, 4).Value = 1.070310489475943
$ws.Cells.Item(25, 5).Value = 1.066589894645995
$ws.Cells.Item(25, 6).Value = 1.07804126149346
$ws.Cells.Item(25, 9).Value = 1.043190084376544
$ws.Cells.Item(25, 10).Value = 1.065928146286071
$ws.Cells.Item(25, 11).Value = 1.073258370636314
$ws.Cells.Item(25, 12).Value = 1.069548703366437
$ws.Cells.Item(25, 13).Value = 1.080966699173843
$ws.Cells.Item(25, 14).Value = 1.025590176298226

Write-Output "applied 380 kV case updates"
